$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-17 down to 13-18
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the K_RISIKO entry
$ws.Range("A12").Value = "K_RISIKO"
$ws.Range("B12").Value = "Risikolage"
$ws.Range("C12").Value = "XXXRisikolage"

# The inserted row already picks up font/alignment/wrap formatting from the
# row it displaced, matching the other data rows - only the border color
# needs to be (re)applied to match the surrounding cells exactly.
$ws.Range("A12:C12").Borders.Color = 15132391
